# Generate Report for Handback
# - Overview sheet: mark zh-cn / de-de status as "Handed back: in sync with en-US"
# - zh-cn sheet: record the handback target file, handback xliff file and handback datetime
# - de-de sheet: record the handback target file, handback xliff file and handback datetime

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet - update the per-language status cells
# ---------------------------------------------------------------------------
$ws_overview.Range("E2").Value = $newStatus
$ws_overview.Range("F2").Value = $newStatus
$ws_overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$ws_overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet - status + handback info
# ---------------------------------------------------------------------------
$ws_zhcn.Range("C2").Value = $newStatus
$ws_zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668

$mdFile = "b202a128-00a6-4338-b4dd-6558c2727ad7.md"
$mdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/12368efc0a8ba2e8a0d199e0b01167a108fcbb95/e2e/b202a128-00a6-4338-b4dd-6558c2727ad7.md"

# Latest Target File (I2) - hyperlink to the handed-back source markdown file
$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("I2"), $mdUrl, "", "", $mdFile) | Out-Null

# Latest Handback File (J2)
$ws_zhcn.Range("J2").Value = "b202a128-00a6-4338-b4dd-6558c2727ad7.d7258a5d7e2df51b1dbc0536218c8a51613f2bf9.zh-cn.xlf"

# Latest Handback DateTime (K2)
$ws_zhcn.Range("K2").Value = "2016-10-21 04:31:55"

$ws_zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws_zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet - status + handback info
# ---------------------------------------------------------------------------
$ws_dede.Range("C2").Value = $newStatus
$ws_dede.Columns.Item(3).ColumnWidth = 29.166666666666668

# Latest Target File (I2) - hyperlink to the handed-back source markdown file
$ws_dede.Hyperlinks.Add($ws_dede.Range("I2"), $mdUrl, "", "", $mdFile) | Out-Null

# Latest Handback File (J2)
$ws_dede.Range("J2").Value = "b202a128-00a6-4338-b4dd-6558c2727ad7.d7258a5d7e2df51b1dbc0536218c8a51613f2bf9.de-de.xlf"

# Latest Handback DateTime (K2)
$ws_dede.Range("K2").Value = "2016-10-21 04:32:13"

$ws_dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws_dede.Columns.Item(10).ColumnWidth = 39.166666666666664
